$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Sema6d"
$ws.Cells.Item(2, 3).Value = "Trem2"
$ws.Cells.Item(2, 4).Value = "M2"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 46.96651666666666
$ws.Cells.Item(2, 8).Value = 140.89955
$ws.Cells.Item(2, 9).Value = 0.5808027674561179
$ws.Cells.Item(2, 10).Value = 0.5808027674561179
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 94.37284199999999
$ws.Cells.Item(2, 14).Value = 283.118526
$ws.Cells.Item(2, 15).Value = 0.9889849826815819
$ws.Cells.Item(2, 16).Value = 0.988984982681582
$ws.Cells.Item(2, 17).Value = 4432.363656673699
$ws.Cells.Item(2, 18).Value = 39891.27291006329
$ws.Cells.Item(2, 19).Value = 0.5744052149140036
$ws.Cells.Item(2, 20).Value = 0.5744052149140036

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Sema6d"
$ws.Cells.Item(3, 3).Value = "Trem2"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 46.96651666666666
$ws.Cells.Item(3, 8).Value = 140.89955
$ws.Cells.Item(3, 9).Value = 0.5808027674561179
$ws.Cells.Item(3, 10).Value = 0.5808027674561179
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 1.051096333333333
$ws.Cells.Item(3, 14).Value = 3.153289
$ws.Cells.Item(3, 15).Value = 0.01101501731841816
$ws.Cells.Item(3, 16).Value = 0.01101501731841816
$ws.Cells.Item(3, 17).Value = 49.36633345777221
$ws.Cells.Item(3, 18).Value = 444.2970011199499
$ws.Cells.Item(3, 19).Value = 0.006397552542114336
$ws.Cells.Item(3, 20).Value = 0.006397552542114336

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Sema6d"
$ws.Cells.Item(4, 3).Value = "Trem2"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 14.34807866666667
$ws.Cells.Item(4, 8).Value = 43.04423600000001
$ws.Cells.Item(4, 9).Value = 0.1774328689611448
$ws.Cells.Item(4, 10).Value = 0.1774328689611448
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 94.37284199999999
$ws.Cells.Item(4, 14).Value = 283.118526
$ws.Cells.Item(4, 15).Value = 0.9889849826815819
$ws.Cells.Item(4, 16).Value = 0.988984982681582
$ws.Cells.Item(4, 17).Value = 1354.068961012904
$ws.Cells.Item(4, 18).Value = 12186.62064911614
$ws.Cells.Item(4, 19).Value = 0.1754784428366812
$ws.Cells.Item(4, 20).Value = 0.1754784428366812

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Sema6d"
$ws.Cells.Item(5, 3).Value = "Trem2"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 14.34807866666667
$ws.Cells.Item(5, 8).Value = 43.04423600000001
$ws.Cells.Item(5, 9).Value = 0.1774328689611448
$ws.Cells.Item(5, 10).Value = 0.1774328689611448
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 1.051096333333333
$ws.Cells.Item(5, 14).Value = 3.153289
$ws.Cells.Item(5, 15).Value = 0.01101501731841816
$ws.Cells.Item(5, 16).Value = 0.01101501731841816
$ws.Cells.Item(5, 17).Value = 15.08121287691156
$ws.Cells.Item(5, 18).Value = 135.730915892204
$ws.Cells.Item(5, 19).Value = 0.00195442612446363
$ws.Cells.Item(5, 20).Value = 0.00195442612446363

# Row 6
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Sema6d"
$ws.Cells.Item(6, 3).Value = "Trem2"
$ws.Cells.Item(6, 4).Value = "M2"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.746361333333333
$ws.Cells.Item(6, 8).Value = 5.239084
$ws.Cells.Item(6, 9).Value = 0.02159605538935411
$ws.Cells.Item(6, 10).Value = 0.02159605538935411
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 94.37284199999999
$ws.Cells.Item(6, 14).Value = 283.118526
$ws.Cells.Item(6, 15).Value = 0.9889849826815819
$ws.Cells.Item(6, 16).Value = 0.988984982681582
$ws.Cells.Item(6, 17).Value = 164.809082185576
$ws.Cells.Item(6, 18).Value = 1483.281739670184
$ws.Cells.Item(6, 19).Value = 0.02135817446523086
$ws.Cells.Item(6, 20).Value = 0.02135817446523086

# Row 7
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Sema6d"
$ws.Cells.Item(7, 3).Value = "Trem2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.746361333333333
$ws.Cells.Item(7, 8).Value = 5.239084
$ws.Cells.Item(7, 9).Value = 0.02159605538935411
$ws.Cells.Item(7, 10).Value = 0.02159605538935411
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 1.051096333333333
$ws.Cells.Item(7, 14).Value = 3.153289
$ws.Cells.Item(7, 15).Value = 0.01101501731841816
$ws.Cells.Item(7, 16).Value = 0.01101501731841816
$ws.Cells.Item(7, 17).Value = 1.835593994141778
$ws.Cells.Item(7, 18).Value = 16.520345947276
$ws.Cells.Item(7, 19).Value = 0.0002378809241232534
$ws.Cells.Item(7, 20).Value = 0.0002378809241232534

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Sema6d"
$ws.Cells.Item(8, 3).Value = "Trem2"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 17.80387266666667
$ws.Cells.Item(8, 8).Value = 53.411618
$ws.Cells.Item(8, 9).Value = 0.2201683081933832
$ws.Cells.Item(8, 10).Value = 0.2201683081933832
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 94.37284199999999
$ws.Cells.Item(8, 14).Value = 283.118526
$ws.Cells.Item(8, 15).Value = 0.9889849826815819
$ws.Cells.Item(8, 16).Value = 0.988984982681582
$ws.Cells.Item(8, 17).Value = 1680.202062159452
$ws.Cells.Item(8, 18).Value = 15121.81855943507
$ws.Cells.Item(8, 19).Value = 0.2177431504656663
$ws.Cells.Item(8, 20).Value = 0.2177431504656663

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Sema6d"
$ws.Cells.Item(9, 3).Value = "Trem2"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 17.80387266666667
$ws.Cells.Item(9, 8).Value = 53.411618
$ws.Cells.Item(9, 9).Value = 0.2201683081933832
$ws.Cells.Item(9, 10).Value = 0.2201683081933832
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 1.051096333333333
$ws.Cells.Item(9, 14).Value = 3.153289
$ws.Cells.Item(9, 15).Value = 0.01101501731841816
$ws.Cells.Item(9, 16).Value = 0.01101501731841816
$ws.Cells.Item(9, 17).Value = 18.71358527906689
$ws.Cells.Item(9, 18).Value = 168.422267511602
$ws.Cells.Item(9, 19).Value = 0.002425157727716944
$ws.Cells.Item(9, 20).Value = 0.002425157727716944
